# Update the column headers on Sheet1:
#   E1: PREREQ_COURSES -> Prerequisites
#   F1: COREQ_COURSES  -> CoRequisites
# (All other cell data is unchanged - only the header labels for the
#  prerequisite/corequisite columns were renamed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Restore the selection to just the header row, matching the saved file.
$ws.Range("A1:I1").Select()
